# Some input files have a worksheet name that is out of the expected
# pattern. Verify the active sheet's name and, if it does not match the
# expected "Sheet1" convention, correct it before any further processing
# (e.g. appending data) takes place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

if ($ws.Name -ne "Sheet1") {
    $ws.Name = "Sheet1"
}
